# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing columns and filling rows 2-4 with the value 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the last existing header cell (G1) onto the new
# header cell (H1) so it picks up the same bold/centered/bordered style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
